# Update the "Pais" worksheet: refresh covid stats and update timestamp,
# matching the commit "Update countries & provincias Spain".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- 1. Update the "last refreshed" timestamp in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 16 de Junio de 2020 a las 06:53"

# --- 2. Refresh per-country statistics (columns: B=Casos totales, ---
# --- C=Nuevos casos, D=Casos activos, E=Recuperados, F=Casos criticos, ---
# --- G=Muertes hoy, H=Muertes) ---

# Row 7: India
$ws.Cells.Item(7, 2).Value = 343091
$ws.Cells.Item(7, 3).Value = 65
$ws.Cells.Item(7, 5).Value = 152856

# Row 18: Pakistan
$ws.Cells.Item(18, 2).Value = 148921
$ws.Cells.Item(18, 3).Value = 4443
$ws.Cells.Item(18, 4).Value = 56390
$ws.Cells.Item(18, 5).Value = 89692
$ws.Cells.Item(18, 7).Value = 110
$ws.Cells.Item(18, 8).Value = 2839

# Row 56: Kazajistan
$ws.Cells.Item(56, 2).Value = 15192
$ws.Cells.Item(56, 3).Value = 383
$ws.Cells.Item(56, 4).Value = 9388
$ws.Cells.Item(56, 5).Value = 5723

# Row 87: El Salvador
$ws.Cells.Item(87, 4).Value = 1927
$ws.Cells.Item(87, 5).Value = 1823
$ws.Cells.Item(87, 7).Value = 2
$ws.Cells.Item(87, 8).Value = 76

# Row 96: Kirguistan
$ws.Cells.Item(96, 2).Value = 2472
$ws.Cells.Item(96, 3).Value = 100
$ws.Cells.Item(96, 4).Value = 1847
$ws.Cells.Item(96, 5).Value = 597
$ws.Cells.Item(96, 7).Value = 1
$ws.Cells.Item(96, 8).Value = 28

# Row 163: Mongolia
$ws.Cells.Item(163, 4).Value = 109
$ws.Cells.Item(163, 5).Value = 88

# Row 184: Butan
$ws.Cells.Item(184, 2).Value = 67
$ws.Cells.Item(184, 3).Value = 1
$ws.Cells.Item(184, 5).Value = 45

# --- 3. Swap rows for two pairs of small territories (name + stats move together) ---
# Santa Sede (row 208) <-> Islas Turcas y Caicos (row 209)
$ws.Cells.Item(208, 1).Value = "Islas Turcas y Caicos"
$ws.Cells.Item(208, 2).Value = 12
$ws.Cells.Item(208, 3).Value = 0
$ws.Cells.Item(208, 4).Value = 11
$ws.Cells.Item(208, 5).Value = 0
$ws.Cells.Item(208, 6).Value = 0
$ws.Cells.Item(208, 7).Value = 0
$ws.Cells.Item(208, 8).Value = 1

$ws.Cells.Item(209, 1).Value = "Santa Sede"
$ws.Cells.Item(209, 2).Value = 12
$ws.Cells.Item(209, 3).Value = 0
$ws.Cells.Item(209, 4).Value = 12
$ws.Cells.Item(209, 5).Value = 0
$ws.Cells.Item(209, 6).Value = 0
$ws.Cells.Item(209, 7).Value = 0
$ws.Cells.Item(209, 8).Value = 0

# Seychelles (row 210) <-> Montserrat (row 211)
$ws.Cells.Item(210, 1).Value = "Montserrat"
$ws.Cells.Item(210, 2).Value = 11
$ws.Cells.Item(210, 3).Value = 0
$ws.Cells.Item(210, 4).Value = 10
$ws.Cells.Item(210, 5).Value = 0
$ws.Cells.Item(210, 6).Value = 0
$ws.Cells.Item(210, 7).Value = 0
$ws.Cells.Item(210, 8).Value = 1

$ws.Cells.Item(211, 1).Value = "Seychelles"
$ws.Cells.Item(211, 2).Value = 11
$ws.Cells.Item(211, 3).Value = 0
$ws.Cells.Item(211, 4).Value = 11
$ws.Cells.Item(211, 5).Value = 0
$ws.Cells.Item(211, 6).Value = 0
$ws.Cells.Item(211, 7).Value = 0
$ws.Cells.Item(211, 8).Value = 0
